# Automatische test-sync: 2025-06-19 21:24:50
# Adds a new mail-log row (row 15) to the "Logs" sheet and bumps the
# matching "Productinformatie" tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append new row 15 ----
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A15").Value = "Is product X op voorraad?"
$logs.Range("B15").Value = "mailmind.test@zohomail.eu"
$logs.Range("C15").Value = "Ik ben geïnteresseerd in product X. Is dit momenteel op voorraad?"
$logs.Range("D15").Value = "Productinformatie"
$logs.Range("F15").Value = "2025-06-19 21:24:10"
$logs.Range("G15").Value = "Nee"

# Extend the conditional-formatting ranges so row 15 is covered too,
# same as the existing rules applied to D2:D14 and G2:G14.
$catRules = $logs.Range("D2:D14").FormatConditions
$catRules.Item(1).ModifyAppliesToRange($logs.Range("D2:D15"))

$answeredRules = $logs.Range("G2:G14").FormatConditions
$answeredRules.Item(1).ModifyAppliesToRange($logs.Range("G2:G15"))

# ---- Dashboard sheet: bump the "Productinformatie" count ----
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B3").Value = 4
